# Auto-generated edit script applying scheduled-runner market data updates
# to Garuda_Profits workbook (per-sheet Leve profit recalculations).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 455.77777
$ws.Range("J32").Value = 455.77777
$ws.Range("L32").Value = 455.77777
$ws.Range("N32").Value = -1107.77777
$ws.Range("H107").Value = 549.2
$ws.Range("I107").Value = 549.2
$ws.Range("K107").Value = 549.2
$ws.Range("M107").Value = 1370.8

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 12346607
$ws.Range("I45").Value = 16667459
$ws.Range("J45").Value = 1313.4286
$ws.Range("K45").Value = 16667459
$ws.Range("L45").Value = 1313.4286
$ws.Range("M45").Value = -16667082
$ws.Range("N45").Value = -2067.4286
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H135").Value = 32865.6
$ws.Range("J135").Value = 32865.6
$ws.Range("L135").Value = 32865.6
$ws.Range("N135").Value = -43005.6

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1633.5264
$ws.Range("I86").Value = 1545.7142
$ws.Range("J86").Value = 1879.4
$ws.Range("K86").Value = 1545.7142
$ws.Range("L86").Value = 1879.4
$ws.Range("M86").Value = -422.7141999999999
$ws.Range("N86").Value = -4125.4
$ws.Range("H89").Value = 1633.5264
$ws.Range("I89").Value = 1545.7142
$ws.Range("J89").Value = 1879.4
$ws.Range("K89").Value = 7728.571
$ws.Range("L89").Value = 9397
$ws.Range("M89").Value = -2112.571
$ws.Range("N89").Value = -20629
$ws.Range("H99").Value = 4000
$ws.Range("I99").Value = 4000
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 4000
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -2502
$ws.Range("N99").ClearContents()
$ws.Range("H107").Value = 1721.3
$ws.Range("I107").Value = 1983.3334
$ws.Range("J107").Value = 1328.25
$ws.Range("K107").Value = 1983.3334
$ws.Range("L107").Value = 1328.25
$ws.Range("M107").Value = -63.33339999999998
$ws.Range("N107").Value = -5168.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 932.35297
$ws.Range("I16").Value = 885
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 885
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -598
$ws.Range("N16").Value = -1574
$ws.Range("H31").Value = 2827919.8
$ws.Range("J31").Value = 6176865
$ws.Range("L31").Value = 6176865
$ws.Range("N31").Value = -6177455
$ws.Range("H34").Value = 2827919.8
$ws.Range("J34").Value = 6176865
$ws.Range("L34").Value = 6176865
$ws.Range("N34").Value = -6177269
$ws.Range("H99").Value = 2848.5386
$ws.Range("I99").Value = 3075
$ws.Range("J99").Value = 2747.889
$ws.Range("K99").Value = 3075
$ws.Range("L99").Value = 2747.889
$ws.Range("M99").Value = -1577
$ws.Range("N99").Value = -5743.889
$ws.Range("H113").Value = 932.35297
$ws.Range("I113").Value = 885
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 885
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 1285
$ws.Range("N113").Value = -5340
$ws.Range("H126").Value = 2848.5386
$ws.Range("I126").Value = 3075
$ws.Range("J126").Value = 2747.889
$ws.Range("K126").Value = 9225
$ws.Range("L126").Value = 8243.667000000001
$ws.Range("M126").Value = -6755
$ws.Range("N126").Value = -13183.667
$ws.Range("H134").Value = 869.4375
$ws.Range("I134").Value = 860.73334
$ws.Range("K134").Value = 2582.20002
$ws.Range("M134").Value = -47.20002000000022

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 684
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 684
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 2052
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value = -5048
$ws.Range("H107").Value = 5101
$ws.Range("I107").Value = 1060
$ws.Range("J107").Value = 5858.6875
$ws.Range("K107").Value = 3180
$ws.Range("L107").Value = 17576.0625
$ws.Range("M107").Value = -1260
$ws.Range("N107").Value = -21416.0625
$ws.Range("H131").Value = 785.77
$ws.Range("J131").Value = 793.30206
$ws.Range("L131").Value = 2379.90618
$ws.Range("N131").Value = -12459.90618

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1364.7273
$ws.Range("I102").Value = 1287.4286
$ws.Range("J102").Value = 1500
$ws.Range("K102").Value = 1287.4286
$ws.Range("L102").Value = 1500
$ws.Range("M102").Value = 334.5714
$ws.Range("N102").Value = -4744
$ws.Range("H126").Value = 2693.625
$ws.Range("I126").Value = 3442.5715
$ws.Range("J126").Value = 2111.111
$ws.Range("K126").Value = 10327.7145
$ws.Range("L126").Value = 6333.333
$ws.Range("M126").Value = -7857.7145
$ws.Range("N126").Value = -11273.333

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 55557468
$ws.Range("I7").Value = 1300
$ws.Range("J7").Value = 100002400
$ws.Range("K7").Value = 1300
$ws.Range("L7").Value = 100002400
$ws.Range("M7").Value = -1188
$ws.Range("N7").Value = -100002624
$ws.Range("H40").Value = 1298
$ws.Range("I40").Value = 1298
$ws.Range("K40").Value = 1298
$ws.Range("M40").Value = -1162
$ws.Range("H46").Value = 1474.1875
$ws.Range("J46").Value = 1676.0769
$ws.Range("L46").Value = 1676.0769
$ws.Range("N46").Value = -2052.0769
$ws.Range("H93").Value = 902005.25
$ws.Range("I93").Value = 1081998.4
$ws.Range("J93").Value = 2039.8
$ws.Range("K93").Value = 1081998.4
$ws.Range("L93").Value = 2039.8
$ws.Range("M93").Value = -1080750.4
$ws.Range("N93").Value = -4535.8
$ws.Range("H126").Value = 55557468
$ws.Range("I126").Value = 1300
$ws.Range("J126").Value = 100002400
$ws.Range("K126").Value = 3900
$ws.Range("L126").Value = 300007200
$ws.Range("M126").Value = -1430
$ws.Range("N126").Value = -300012140

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1305.7368
$ws.Range("I122").Value = 883.7727
$ws.Range("J122").Value = 1885.9375
$ws.Range("K122").Value = 2651.3181
$ws.Range("L122").Value = 5657.8125
$ws.Range("M122").Value = -201.3181
$ws.Range("N122").Value = -10557.8125
$ws.Range("H126").Value = 35719216
$ws.Range("J126").Value = 1222.1666
$ws.Range("L126").Value = 3666.4998
$ws.Range("N126").Value = -8606.4998

